# Insert a new price-record row for "Vega Monumental Concepción" / Mango at
# row 91 (pushing the existing rows 91-194 down to 92-195), and populate it
# with a new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 91..194 down to 92..195, leaving row 91 free for the
# new record (mirrors the other rows' constant columns A,B,C,E-L,Q,T).
$ws.Rows.Item(91).Insert()

$ws.Cells.Item(91, 1).Value2 = 11
$ws.Cells.Item(91, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(91, 3).Value2 = "Bíobío"
$ws.Cells.Item(91, 4).Value2 = 45195
$ws.Cells.Item(91, 5).Value2 = 8
$ws.Cells.Item(91, 6).Value2 = "Fruta"
$ws.Cells.Item(91, 7).Value2 = 100108
$ws.Cells.Item(91, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(91, 9).Value2 = 100108002
$ws.Cells.Item(91, 10).Value2 = "Mango"
$ws.Cells.Item(91, 11).Value2 = "Sin especificar"
$ws.Cells.Item(91, 12).Value2 = "Primera"
$ws.Cells.Item(91, 13).Value2 = 100
$ws.Cells.Item(91, 14).Value2 = 11000
$ws.Cells.Item(91, 15).Value2 = 11000
$ws.Cells.Item(91, 16).Value2 = 11000
$ws.Cells.Item(91, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(91, 18).Value2 = "Brasil"
$ws.Cells.Item(91, 19).Value2 = 2750
$ws.Cells.Item(91, 20).Value2 = 4
